# A previous merge failed, so manually re-add the missing "destinations" data
# update: the "Latest period (release date)" for KS4/KS5 destinations rows
# moves from the 24/10/24 publication date to the 27/02/25 one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = "Aug 2022 -  Jul 2023 (21/22 learners) (27/02/25)"

$ws.Range("C11").Value = $newDate
$ws.Range("C12").Value = $newDate

$ws.Range("C13").Select()
